$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new 2035 IPA runs (rows 32 and 33) ---
# Duplicate the formatting of the existing 2035 IPA rows (29-31) onto the
# two new rows by copying a full-row range (values+formats), then
# overwrite the cells that actually differ for the new runs.

$ws.Range("A30:S30").Copy($ws.Range("A32:S32"))
$ws.Range("A30:S30").Copy($ws.Range("A33:S33"))

# Row 32: 2035_TM160_IPA_04 - WFH rate at 25%, with trn hesistance on
$ws.Range("A32").Value = 2035
$ws.Range("B32").Value = "2035_TM160_IPA_04"
$ws.Range("C32").Value = "RTP2025_IP"
$ws.Range("D32").Value = "Future year"
$ws.Range("E32").Value = "WFH rate at 25%, with trn hesistance on"
$ws.Range("F32").Value = "petrale"
$ws.Range("G32").Value = "n/a"
$ws.Range("H32").Value = "current"
$ws.Range("I32").Value = "M:\Application\Model One\RTP2021\Blueprint\INPUT_DEVELOPMENT\Networks\BlueprintNetworks_64\net_2035_Blueprint_tollscsv"
$ws.Range("J32").Value = "model3-a"
$ws.Range("K32").Value = "https://app.asana.com/0/1204085012544660/1205561944199029/f"
$ws.Range("L32").Value = 18.64
$ws.Range("M32").Value = "na"
$ws.Range("N32").Value = "na"
$ws.Range("O32").Value = 0.95
$ws.Range("P32").Value = 0.86
$ws.Range("Q32").Value = 120
$ws.Range("R32").Value = 0
$ws.Range("S32").Value = 45

# Row 33: 2035_TM160_IPA_05 - WFH rate at 25%, with trn hesistance off
$ws.Range("A33").Value = 2035
$ws.Range("B33").Value = "2035_TM160_IPA_05"
$ws.Range("C33").Value = "RTP2025_IP"
$ws.Range("D33").Value = "Future year"
$ws.Range("E33").Value = "WFH rate at 25%, with trn hesistance off"
$ws.Range("F33").Value = "petrale"
$ws.Range("G33").Value = "n/a"
$ws.Range("H33").Value = "current"
$ws.Range("I33").Value = "M:\Application\Model One\RTP2021\Blueprint\INPUT_DEVELOPMENT\Networks\BlueprintNetworks_64\net_2035_Blueprint_tollscsv"
$ws.Range("J33").Value = "model3-b"
$ws.Range("K33").Value = "https://app.asana.com/0/1204085012544660/1205561944199034/f"
$ws.Range("L33").Value = 18.64
$ws.Range("M33").Value = "na"
$ws.Range("N33").Value = "na"
$ws.Range("O33").Value = 0.95
$ws.Range("P33").Value = 0.86
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = 0
$ws.Range("S33").Value = 0

# --- Update the sheet view to reflect scrolling down to the new rows ---
$ws.Range("A33").Select()

Write-Host "Edit complete"
